$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Step Response" - add the "1Tau" analysis columns (E:J) next to the
# existing Before/After step-response readings, with new header labels,
# a handful of newly recorded voltages, and the derived formulas.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Step Response")

# New header labels for columns E:J (entered in this order so new shared
# strings land in the same sequence as the source workbook).
$ws3.Range("E1").Value = "Top Tank 1Tau Height (in)"
$ws3.Range("G1").Value = "Top Tank Tau (s)"
$ws3.Range("J1").Value = "Qi2/H1 Gain"
$ws3.Range("F1").Value = "Bottom Tank ""1Tau"" Height (in)"
$ws3.Range("H1").Value = "Bottom Tank ""Tau"" (s)"
$ws3.Range("I1").Value = "TopHeight/GPM Gain"

# --- Step response block (rows 3-4) ---------------------------------------
$ws3.Range("C3").Value = 7.25
$ws3.Range("C4").Value = 9.25
$ws3.Range("D4").Value = 7.5
$ws3.Range("E4").Formula = "=C3+0.632*(C4-C3)"
$ws3.Range("F4").Formula = "=D3+0.632*(D4-D3)"
$ws3.Range("G4").Value = 88
$ws3.Range("H4").Value = 163
$ws3.Range("I4").Formula = "=(C4-C3)/(((B4*10.399-6.1525)*0.062+0.15)-((B3*10.399-6.1525)*0.062+0.15))"
$ws3.Range("J4").Formula = "=((B4*10.399-6.1525)*0.062+0.15)/C4"

# --- Step response block (rows 6-7) ----------------------------------------
$ws3.Range("C6").Value = 9.25
$ws3.Range("D6").Value = 7.5
$ws3.Range("D7").Value = 9.5
$ws3.Range("E7").Formula = "=C6+0.632*(C7-C6)"
$ws3.Range("F7").Formula = "=D6+0.632*(D7-D6)"
$ws3.Range("G7").Value = 105
$ws3.Range("H7").Value = 169
$ws3.Range("I7").Formula = "=(C7-C6)/(((B7*10.399-6.1525)*0.062+0.15)-((B6*10.399-6.1525)*0.062+0.15))"
$ws3.Range("J7").Formula = "=((B7*10.399-6.1525)*0.062+0.15)/C7"

# --- Step response block (rows 9-10) ----------------------------------------
$ws3.Range("D9").Value = 9.5
$ws3.Range("E10").Formula = "=C9+0.632*(C10-C9)"
$ws3.Range("F10").Formula = "=D9+0.632*(D10-D9)"
$ws3.Range("G10").Value = 109
$ws3.Range("H10").Value = 192
$ws3.Range("I10").Formula = "=(C10-C9)/(((B10*10.399-6.1525)*0.062+0.15)-((B9*10.399-6.1525)*0.062+0.15))"
$ws3.Range("J10").Formula = "=((B10*10.399-6.1525)*0.062+0.15)/C10"

# --- Averages of the Qi2/H1 and TopHeight/GPM gains -------------------------
$ws3.Range("I13").Formula = "=AVERAGE(I4,I7,I10)"
$ws3.Range("J13").Formula = "=AVERAGE(J4,J7,J10)"

# --- Disturbance block (rows 15-16) -----------------------------------------
$ws3.Range("E16").Formula = "=C15+0.632*(C16-C15)"
$ws3.Range("F16").Formula = "=D15+0.632*(D16-D15)"

# ---------------------------------------------------------------------------
# Chart: give the bottom-tank linear-fit chart an explicit title.
# ---------------------------------------------------------------------------
$wsPressure = $wb.Worksheets.Item("Part 2 - Pressure Transducer")
$bottomLinearChart = $wsPressure.ChartObjects().Item(4).Chart
$bottomLinearChart.HasTitle = $true
$bottomLinearChart.ChartTitle.Text = "Bottom Tank Height vs. Vhbttm Linear"

# ---------------------------------------------------------------------------
# Restore each sheet's selection, ending on "Step Response" so it stays the
# active tab.
# ---------------------------------------------------------------------------
$wsValve = $wb.Worksheets.Item("Part 1 - Valve")
$wsValve.Activate()
$wsValve.Range("E6").Select()

$wsPressure.Activate()
$wsPressure.Range("D41").Select()

$ws3.Activate()
$ws3.Range("C6").Select()
